$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet carries no real protection password in the workbook, but the
# runtime defaults newly-opened sheets to a protected state; unprotect so
# the cell writes below are allowed.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer banner
# (shared string used by cell A59) from 2021-05-26 to 2021-05-27.
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

# Refreshed "Weight" (column D) and "Percent Change" (column E) figures
# for each holding row (2-55), plus the portfolio total's percent change
# in E56.
$changes = @(
    @("D2", 0.023052699865487),
    @("E2", -0.01802600472813254),
    @("D3", 0.01770293197876988),
    @("E3", 0.0005937067088859216),
    @("D4", 0.01832537823147021),
    @("E4", 0.002039255671679951),
    @("D5", 0.02025519517816308),
    @("E5", 0.0145963658028001),
    @("D6", 0.01930167291675629),
    @("E6", 0.01016456921587605),
    @("D7", 0.02663342022537693),
    @("E7", 0.04090267983074747),
    @("D8", 0.01907239158915122),
    @("E8", 0.002510460251045954),
    @("D9", 0.0192222613873311),
    @("E9", 0.02559740785743214),
    @("D10", 0.0186897370135387),
    @("E10", 0.02111972006998242),
    @("D11", 0.01942370975241705),
    @("E11", 0.02329776040883802),
    @("D12", 0.0190319072800325),
    @("E12", 0.0158311345646438),
    @("D13", 0.02049479221914936),
    @("E13", -0.0006457862447530882),
    @("D14", 0.01925145680256095),
    @("E14", -0.01245576786978075),
    @("D15", 0.01783878797763943),
    @("E15", -0.009601536245799402),
    @("D16", 0.01804082025102997),
    @("E16", 0.009817671809256634),
    @("D17", 0.01552222909720186),
    @("E17", 0.01241379310344826),
    @("D18", 0.01572192573737401),
    @("E18", 0.08940774487471526),
    @("D19", 0.01663982959220039),
    @("E19", 0.01865671641791056),
    @("D20", 0.01960764086836508),
    @("E20", 0.001389716100853589),
    @("D21", 0.02036905729755948),
    @("E21", 0.003955968352252981),
    @("D22", 0.02097768438705102),
    @("E22", 0.01403798513625087),
    @("D23", 0.0181945827712405),
    @("E23", 0.01171373555840827),
    @("D24", 0.02114915879250099),
    @("E24", -0.002576845205227274),
    @("D25", 0.02156568004978014),
    @("E25", 0.01444043321299637),
    @("D26", 0.02078830346026008),
    @("E26", 0.002078534913768859),
    @("D27", 0.01935694956959146),
    @("E27", 0.005087881591119281),
    @("D28", 0.028004242288469),
    @("E28", -0.006880733944954143),
    @("D29", 0.01931568671606661),
    @("E29", 0.01420798065296247),
    @("D30", 0.01269183090871898),
    @("E30", 0.005183413078150068),
    @("D31", 0.009407735967564286),
    @("E31", -0.0007241129616221142),
    @("D32", 0.01661238590188433),
    @("E32", 0.01459854014598538),
    @("D33", 0.01901614175580838),
    @("E33", 0.005639655684179345),
    @("D34", 0.0182825582891331),
    @("E34", -0.007835455435847183),
    @("D35", 0.02010668783269393),
    @("E35", -0.02722063037249289),
    @("D36", 0.01755812271922984),
    @("E36", 0.0419243986254294),
    @("D37", 0.01919773723853803),
    @("E37", -0.01706308169596682),
    @("D38", 0.01931179399403597),
    @("E38", 0),
    @("D39", 0.02534006333069472),
    @("E39", 0.0131498095121052),
    @("D40", 0.01563550730829366),
    @("E40", 0.02270577105014193),
    @("D41", 0.02156879422740465),
    @("E41", 0.006064106266243208),
    @("D42", 0.01938906452634429),
    @("E42", 0.01621209231356113),
    @("D43", 0.01992723334708114),
    @("E43", 0.0131761442441054),
    @("D44", 0.01749914798046555),
    @("E44", 0.008953696597595151),
    @("D45", 0.02103043077056628),
    @("E45", 0.01457658491439151),
    @("D46", 0.01946166379221585),
    @("E46", 0.002640264026402495),
    @("D47", 0.01823370462764849),
    @("E47", 0.01972651871777642),
    @("D48", 0.01606915654250765),
    @("E48", 0.0159883720930234),
    @("D49", 0.01745885830744836),
    @("E49", 0.04384615384615387),
    @("D50", 0.01713186965687407),
    @("E50", 0.03435582822085892),
    @("D51", 0.01671826794111791),
    @("E51", 0.02291169451073993),
    @("D52", 0.01847485875744703),
    @("E52", 0.02739148756847865),
    @("D53", 0.0155845126496922),
    @("E53", 0.03057324840764308),
    @("D54", 0.007695911454587603),
    @("E54", -0.002655538694992465),
    @("D55", 0.007045826875469684),
    @("E55", -0.007955801104972404),
    @("E56", 0.01095645542745705)
)

foreach ($pair in $changes) {
    $ws.Range($pair[0]).Value = $pair[1]
}
